$wb = $excel.ActiveWorkbook

# The "smoke" sheet becomes the active sheet/tab (workbookView activeTab="1",
# tabSelected moves from AppControl's sheetView to smoke's sheetView).
$smoke = $wb.Worksheets.Item("smoke")
$smoke.Activate()

# Rows 3-18 in column B ("Flag") flip from "Y" to "N".
$smoke.Range("B3:B18").Value = "N"

# Selection / scroll position on the smoke sheet moves to B3:B21 (active
# cell B3), with the view scrolled so row 10 is at the top.
$smoke.Range("B3:B21").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
